$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Starting/Ending SoC (%) values
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 24

# Row 8 - label change only
$ws.Range("A8").Value = "Total distance covered (km)"

# Row 9 - label change only
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

# Row 10 - label + value change
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 75

# Row 12 - label change only
$ws.Range("A12").Value = "Peak Power(kW)"

# Row 13 - label change only
$ws.Range("A13").Value = "Average Power(kW)"

# Row 14 - label change only
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# Row 15 - label + value change (sign flip)
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.05122061903286506

# Row 16 - label + value change (was Lowest Cell Voltage, becomes Highest Cell Voltage)
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.335

# Row 17 - label + value change (was Highest Cell Voltage, becomes Lowest Cell Voltage)
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.021

# Row 18 - label change only
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# Row 19 - label change only
$ws.Range("A19").Value = "Minimum Temperature(C)"

# Row 20 - label change only
$ws.Range("A20").Value = "Maximum Temperature(C)"

# Row 21 - label change + value now populated
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 18

# Row 22 - label change only
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"

# Row 23 - label change only
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"

# Row 24 - label change only
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"

# Row 25 - label change only
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"

# Row 26 - label change only
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"

# Row 27 - label change only
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# Row 28 - label change (was lowest cell temp, becomes highest cell temp)
$ws.Range("A28").Value = "highest cell temp(C)"

# Row 29 - label change (was highest cell temp, becomes lowest cell temp)
$ws.Range("A29").Value = "lowest cell temp(C)"

# Row 30 - label change only
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# Row 31 - label + value change (was Maximum BMS Temperature in C, becomes Battery Voltage(V))
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

# Row 32 - label + value change (was Battery Voltage, becomes Total energy charged(kWh))
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.511638660833333

# Row 33 - label + value change (was Total energy charged in kWh, becomes Electricity consumption units(kW))
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001449929654728105

# Row 34 - label + value change (was Electricity consumption units in kW, becomes Idling time percentage)
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 21.82369285173023

# Row 35 - label + value change (was Idling time percentage, becomes Time spent in 0-10 km/h)
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 6.05792708596447

# Row 36 - label + value change (was Time spent in 0-10 km/h, becomes Time spent in 10-20 km/h)
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 3.279447671971036

# Row 37 - label + value change (was Time spent in 10-20 km/h, becomes Time spent in 20-30 km/h)
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 5.430664309168982

# Row 38 - label + value change (was Time spent in 20-30 km/h, becomes Time spent in 30-40 km/h)
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 10.47823524459039

# Row 39 - label + value change (was Time spent in 30-40 km/h, becomes Time spent in 40-50 km/h)
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 8.832196682663973

# Row 40 - label + value change (was Time spent in 40-50 km/h, becomes Time spent in 50-60 km/h)
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 17.6433442788583

# Row 41 - label + value change (was Time spent in 50-60 km/h, becomes Time spent in 60-70 km/h)
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 14.4270438662962

# Row 42 - label + value change (was Time spent in 60-70 km/h, becomes Time spent in 70-80 km/h)
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 11.9390418455839

# Row 43 - new row
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
